$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "LSL code and sound effects..." -> "Scripting" + bookmark(_GoBack)
#           + " and sound effects..."
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("LSL code", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Scripting", 2)

# ---------------------------------------------------------------------------
# Change 3 (bookmark source removal): delete the existing "_GoBack" bookmark
# that currently sits between "Worked in a fast-paced, on" and
# "line environment..." and merge those two runs back into one run with the
# combined text "Worked in a fast-paced, online environment...".
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$full = $d.Content.Text
$targetText3 = "Worked in a fast-paced, online environment to deliver tutoring services to students nation-wide."
$idx3 = $full.IndexOf("Worked in a fast-paced, on")
$end3 = $idx3 + $targetText3.Length
$r3 = $d.Range($idx3, $end3)
$r3.Text = "PLACEHOLDER_THREE"
$full = $d.Content.Text
$idx3p = $full.IndexOf("PLACEHOLDER_THREE")
$r3b = $d.Range($idx3p, $idx3p + "PLACEHOLDER_THREE".Length)
$r3b.Text = $targetText3

# ---------------------------------------------------------------------------
# Change 1b (bookmark target insertion): add the "_GoBack" bookmark between
# "Scripting" and " and sound effects..." (collapsed range, no selected text).
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx1 = $full.IndexOf("Scripting and sound")
$pos1 = $idx1 + "Scripting".Length
$r1 = $d.Range($pos1, $pos1)
$null = $d.Bookmarks.Add("_GoBack", $r1)

# ---------------------------------------------------------------------------
# Change 2: merge "modern SCM and time tracking utilities" + " " runs into a
# single run "modern SCM and time tracking utilities " (trailing space),
# while leaving the surrounding "Used " / "to ensure..." runs untouched.
# Pin the merge boundaries with throw-away bookmarks so the auto-merge that
# happens on a Range.Text assignment doesn't swallow the neighboring runs.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx2 = $full.IndexOf("modern SCM and time tracking utilities")
$targetText2 = "modern SCM and time tracking utilities "
$end2 = $idx2 + $targetText2.Length

$rBefore = $d.Range($idx2, $idx2)
$null = $d.Bookmarks.Add("TempBefore", $rBefore)
$rAfter = $d.Range($end2, $end2)
$null = $d.Bookmarks.Add("TempAfter", $rAfter)

$r2 = $d.Range($idx2, $end2)
$r2.Text = "PLACEHOLDER_TWO"
$full = $d.Content.Text
$idx2p = $full.IndexOf("PLACEHOLDER_TWO")
$r2b = $d.Range($idx2p, $idx2p + "PLACEHOLDER_TWO".Length)
$r2b.Text = $targetText2

$d.Bookmarks.Item("TempBefore").Delete()
$d.Bookmarks.Item("TempAfter").Delete()

Write-Host "All changes applied"
